$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "EMPID"
$ws.Range("B1").Value = "LNAME"
$ws.Range("C1").Value = "FNAME"
$ws.Range("D1").Value = "BDATE"

$ws.Range("A2:E6").Select() | Out-Null
